$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SkillData")

$ws.Range("B2").Value = "Skill_K_King001"
$ws.Range("B3").Value = "Skill_S_Knight001"
$ws.Range("B4").Value = "Skill_S_Archer001"
$ws.Range("B5").Value = "Skill_S_Magic001"
$ws.Range("B6").Value = "Skill_H_Knight001"
$ws.Range("B7").Value = "Skill_H_Archer001"
$ws.Range("B8").Value = "Skill_H_Magic001"

$ws.Range("G9").Select()
